$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-19 down to 8-20.
$ws.Rows("7:7").Insert()

# --- Fill the new row 7: "Reverse Words in a String III" ---
$ws.Range("A7").Value = "Arrays"
$ws.Range("B7").Value = "Reverse Words in a String III"
$ws.Range("C7").Value = "https://leetcode.com/problems/reverse-words-in-a-string-iii"
$ws.Range("D7").Value = "ReverseWordsInAStringIII_TraverseReverse"
$ws.Range("E7").Value = "Traverse and Reverse each character one by one"

# F7: rich-text note with mixed red/black runs
$noteText = "Their's only one edge case in here, the last word does not have a space after its last character, you need two loop (loop one will iterate whole the string , and loop two will do reversing string), of course you need and if statement to check whether their's space or is lastIndex (if ((strIndex == s.length() - 1) || s.charAt(strIndex) == ' ')) "
$ws.Range("F7").Value = $noteText

$run2 = $ws.Range("F7").Characters(99,17)
$run2.Font.Size = 14
$run2.Font.Name = "Arial (Body)"
$run2.Font.Color = 255

$run3 = $ws.Range("F7").Characters(116,2)
$run3.Font.Size = 14
$run3.Font.Name = "Arial (Body)"
$run3.Font.Color = 0

$run4 = $ws.Range("F7").Characters(118,78)
$run4.Font.Size = 14
$run4.Font.Name = "Arial (Body)"
$run4.Font.Color = 255

$run5 = $ws.Range("F7").Characters(196,85)
$run5.Font.Size = 14
$run5.Font.Name = "Arial (Body)"
$run5.Font.Color = 0

$run6 = $ws.Range("F7").Characters(281,66)
$run6.Font.Size = 14
$run6.Font.Name = "Arial (Body)"
$run6.Font.Color = 255

# --- Old row 7 (now row 8, "Running Sum of 1d Array") gains a Class name ---
$ws.Range("D8").Value = "RunningSumOf1dArray"

# --- Rebuild all hyperlinks at their (possibly shifted) locations ---
# Row-insert does not auto-shift the existing Hyperlinks collection, so clear
# it and re-add every link in the same relative order to land on the correct
# cells with the correct relationship ids.
$ws.Range("C3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "https://leetcode.com/problems/two-sum/")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://leetcode.com/problems/valid-palindrome")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/problems/merge-sorted-array")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://leetcode.com/problems/is-subsequence")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://leetcode.com/problems/reverse-string")
$ws.Hyperlinks.Add($ws.Range("C11"), "https://leetcode.com/problems/squares-of-a-sorted-array/")
$ws.Hyperlinks.Add($ws.Range("C14"), "https://leetcode.com/problems/subarray-product-less-than-k/")
$ws.Hyperlinks.Add($ws.Range("C16"), "https://leetcode.com/problems/maximum-average-subarray-i/")
$ws.Hyperlinks.Add($ws.Range("C17"), "https://leetcode.com/problems/max-consecutive-ones-iii/")
$ws.Hyperlinks.Add($ws.Range("C18"), "https://leetcode.com/problems/number-of-ways-to-split-array/")
$ws.Hyperlinks.Add($ws.Range("C19"), "https://leetcode.com/problems/number-of-ways-to-split-array/")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://leetcode.com/problems/running-sum-of-1d-array/")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://leetcode.com/problems/minimum-value-to-get-positive-step-by-step-sum/")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://leetcode.com/problems/k-radius-subarray-averages/")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/problems/reverse-words-in-a-string-iii")

# Restore the plain "Link" look (style used throughout column C) on every
# hyperlinked cell, since Hyperlinks.Add re-styles cells with the built-in
# Hyperlink font/underline.
$ws.Range("C2").Copy()
$ws.Range("C2:C20").PasteSpecial(-4122)
$ws.Range("C2").Select()
$excel.CutCopyMode = 0

# --- Update the view selection to match the target state ---
$ws.Range("G7").Select()

Write-Output "done"
